$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# Row 18: fill in a date, description, and amount (previously blank placeholder row)
$date = Get-Date -Year 2016 -Month 4 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Range("B18").Value = $date.Date
$ws.Range("C18").Value = "In thử bản 38tr. Kẹp lò xo 15k, 3 tờ A3 12k, 35 tờ A4 màu 35k"
$ws.Range("E18").Value = 62

# Extend the total formula to include the newly filled row 18
$ws.Range("E21").Formula = "=SUM(E4:E18)"

# Move the active selection to C19 (next empty row), matching the updated UI state
$ws.Range("C19").Select()
